# Apply updated crypto price/volume snapshot (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.237.02'
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").Value = '1.902.80'
$ws.Range("E3").Value = '  +0.76%  '

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'308.14"
$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("E7").Value = '  +0.41%  '

$ws.Range("D8").Value = "'0.3773"
$ws.Range("E8").Value = '  +0.51%  '

$ws.Range("D9").Value = "'0.07279"
$ws.Range("E9").Value = '  +1.09%  '

$ws.Range("D10").Value = "'21.19"
$ws.Range("E10").Value = '  +0.34%  '

$ws.Range("E11").Value = '  +0.29%  '

$ws.Range("D12").Value = "'0.08277"
$ws.Range("E12").Value = '  +8.37%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.918.85'
$ws.Range("E13").Value = '  +1.61%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = "'96.69"
$ws.Range("E14").Value = '  +2.45%  '

$ws.Range("D15").Value = "'5.281"
$ws.Range("E15").Value = '  +0.80%  '

$ws.Range("D16").Value = "'1.0000"
$ws.Range("E16").Value = '  -0.05%  '

$ws.Range("D17").Value = "'0.000008639"
$ws.Range("E17").Value = '  +1.66%  '

$ws.Range("D18").Value = "'14.57"
$ws.Range("E18").Value = '  +0.77%  '

$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = '  +0.06%  '

$ws.Range("D20").Value = '27.257.33'
$ws.Range("E20").Value = '  +0.71%  '

$ws.Range("D21").Value = "'5.088"
$ws.Range("E21").Value = '  +0.77%  '

$ws.Range("D22").Value = '2.149.14'
$ws.Range("E22").Value = '  +1.35%  '

$ws.Range("D23").Value = "'10.66"
$ws.Range("E23").Value = '  +0.70%  '

$ws.Range("E24").Value = '  +0.75%  '

$ws.Range("D25").Value = "'2.328"
$ws.Range("E25").Value = '  +1.23%  '

$ws.Range("D26").Value = "'147.00"
$ws.Range("E26").Value = '  +0.78%  '

$ws.Range("D27").Value = "'1.746"
$ws.Range("E27").Value = '  +0.77%  '

$ws.Range("E28").Value = '  +0.81%  '

$ws.Range("D29").Value = "'115.19"
$ws.Range("E29").Value = '  +0.85%  '

$ws.Range("D30").Value = "'4.836"
$ws.Range("E30").Value = '  +0.98%  '

$ws.Range("E31").Value = '  -0.35%  '

$ws.Range("D32").Value = "'0.09245"
$ws.Range("E32").Value = '  +0.59%  '

$ws.Range("D33").Value = "'0.05078"
$ws.Range("E33").Value = '  +0.86%  '

$ws.Range("D34").Value = "'0.7974"
$ws.Range("E34").Value = '  +3.85%  '

$ws.Range("D35").Value = "'1.239"
$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("D36").Value = "'3.425"
$ws.Range("E36").Value = '  +4.67%  '

$ws.Range("D37").Value = "'2.962"
$ws.Range("E37").Value = '  -0.50%  '

$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("D39").Value = "'0.5712"
$ws.Range("E39").Value = '  +2.04%  '

$ws.Range("D40").Value = "'0.02006"
$ws.Range("E40").Value = '  +0.88%  '

$ws.Range("D41").Value = "'1.078"
$ws.Range("E41").Value = '  +0.45%  '

$ws.Range("D42").Value = "'9.031"
$ws.Range("E42").Value = '  +0.26%  '

$ws.Range("D43").Value = "'6.593"
$ws.Range("E43").Value = '  -0.34%  '

$ws.Range("D44").Value = "'116.82"
$ws.Range("E44").Value = '  -1.72%  '

$ws.Range("D45").Value = "'0.1518"
$ws.Range("E45").Value = '  +0.86%  '

$ws.Range("D46").Value = "'0.4858"
$ws.Range("E46").Value = '  +0.70%  '

$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = '  +0.05%  '

$ws.Range("D48").Value = "'10.15"
$ws.Range("E48").Value = '  -0.26%  '

$ws.Range("D49").Value = "'1.627"

$ws.Range("D50").Value = "'37.68"
$ws.Range("E50").Value = '  +0.05%  '

$ws.Range("D51").Value = "'63.93"
$ws.Range("E51").Value = '  +0.01%  '
